# Apply updated cryptocurrency ranking data (price & 1h volume change)
# as scraped on Tue Feb 20 22:47:08 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "52.268.41"
$ws.Cells.Item(2, 5).Value = "  +0.52%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.993.14"
$ws.Cells.Item(3, 5).Value = "  +0.78%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.05%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'354.08"
$ws.Cells.Item(5, 5).Value = "  +0.16%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'108.42"
$ws.Cells.Item(6, 5).Value = "  -3.45%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.563"
$ws.Cells.Item(7, 5).Value = "  +0.11%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.01%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.625"
$ws.Cells.Item(9, 5).Value = "  -0.94%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'38.43"
$ws.Cells.Item(10, 5).Value = "  -3.39%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +2.49%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.0861"
$ws.Cells.Item(12, 5).Value = "  -4.35%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'19.34"
$ws.Cells.Item(13, 5).Value = "  -3.29%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "3.464.75"
$ws.Cells.Item(14, 5).Value = "  +0.71%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'7.78"
$ws.Cells.Item(15, 5).Value = "  -2.18%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "3.000.70"
$ws.Cells.Item(16, 5).Value = "  +0.64%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +2.60%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "52.265.62"
$ws.Cells.Item(18, 5).Value = "  +0.30%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'3.52"
$ws.Cells.Item(19, 5).Value = "  +5.92%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'7.56"
$ws.Cells.Item(20, 5).Value = "  -2.37%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'13.68"
$ws.Cells.Item(21, 5).Value = "  -5.83%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0976"
$ws.Cells.Item(22, 5).Value = "  -1.84%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'69.61"
$ws.Cells.Item(23, 5).Value = "  -2.58%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'264.49"
$ws.Cells.Item(24, 5).Value = "  -2.44%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -2.34%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'0.180"
$ws.Cells.Item(26, 5).Value = "  -0.67%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'26.91"
$ws.Cells.Item(27, 5).Value = "  -2.37%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'7.53"
$ws.Cells.Item(28, 5).Value = "  -1.84%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.05%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -1.80%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "RenderToken"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(31, 4).Value = "'6.48"
$ws.Cells.Item(31, 5).Value = "  +4.04%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Cosmos"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(32, 4).Value = "'10.36"
$ws.Cells.Item(32, 5).Value = "  -3.69%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'36.55"
$ws.Cells.Item(33, 5).Value = "  -3.40%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +10.69%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'50.94"
$ws.Cells.Item(35, 5).Value = "  -4.09%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -1.10%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -0.05%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'3.22"
$ws.Cells.Item(38, 5).Value = "  -6.87%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'17.95"
$ws.Cells.Item(39, 5).Value = "  -5.49%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'2.71"
$ws.Cells.Item(41, 5).Value = "  +0.46%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.19%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'22.86"
$ws.Cells.Item(43, 5).Value = "  -4.68%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'123.73"
$ws.Cells.Item(44, 5).Value = "  +9.08%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -1.43%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "2.125.82"
$ws.Cells.Item(46, 5).Value = "  -2.77%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'3.40"
$ws.Cells.Item(47, 5).Value = "  -4.77%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'2.41"
$ws.Cells.Item(48, 5).Value = "  -5.12%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.250"
$ws.Cells.Item(49, 5).Value = "  +2.27%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "RocketPoolETH"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(50, 4).Value = "3.285.51"
$ws.Cells.Item(50, 5).Value = "  +0.57%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "BEAM"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Cells.Item(51, 4).Value = "'0.0332"
$ws.Cells.Item(51, 5).Value = "  -2.64%  "

